$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-04-17 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-04-18 Thursday", 2) | Out-Null
$d.Content.Find.Execute("261÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "804÷8=", 2) | Out-Null
$d.Content.Find.Execute("225÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "684÷6=", 2) | Out-Null
$d.Content.Find.Execute("110÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "655÷3=", 2) | Out-Null
$d.Content.Find.Execute("186÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "783÷3=", 2) | Out-Null
$d.Content.Find.Execute("902÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "295÷8=", 2) | Out-Null
$d.Content.Find.Execute("299÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "433÷8=", 2) | Out-Null
$d.Content.Find.Execute("788÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "629÷9=", 2) | Out-Null
$d.Content.Find.Execute("895÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "194÷9=", 2) | Out-Null
$d.Content.Find.Execute("989÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "251÷7=", 2) | Out-Null
$d.Content.Find.Execute("933÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "705÷7=", 2) | Out-Null
$d.Content.Find.Execute("586÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "453÷6=", 2) | Out-Null
$d.Content.Find.Execute("424÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "971÷3=", 2) | Out-Null
$d.Content.Find.Execute("850÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "114÷4=", 2) | Out-Null
$d.Content.Find.Execute("917÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "373÷4=", 2) | Out-Null
$d.Content.Find.Execute("490÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "839÷6=", 2) | Out-Null
$d.Content.Find.Execute("475÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "468÷4=", 2) | Out-Null
$d.Content.Find.Execute("642÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "491÷2=", 2) | Out-Null
$d.Content.Find.Execute("787÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "249÷3=", 2) | Out-Null
$d.Content.Find.Execute("427÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "312÷8=", 2) | Out-Null
$d.Content.Find.Execute("260÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "316÷2=", 2) | Out-Null
$d.Content.Find.Execute("643÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "239÷4=", 2) | Out-Null
$d.Content.Find.Execute("491÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "420÷4=", 2) | Out-Null
$d.Content.Find.Execute("831÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "707÷3=", 2) | Out-Null
$d.Content.Find.Execute("722÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "641÷5=", 2) | Out-Null
$d.Content.Find.Execute("878÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "914÷6=", 2) | Out-Null
